$d = $word.ActiveDocument

# 1) Representative's name was corrected/expanded:
#    "Manuel Dias" -> "Manuel Inácio Veladas Dias"
$d.Content.Find.Execute("Manuel Dias", $true, $true, $false, $false, $false,
                         $true, 1, $false, "Manuel Inácio Veladas Dias", 2)

# 2) Typo fix: the run held "programas" immediately followed by a run starting
#    with "s, projetos ..." (rendering as the misspelling "programass, projetos").
#    Trim the run's text to "programa" so the combined text reads "programas, projetos".
#    (MatchWholeWord must stay False: "programas" here abuts a following "s" with
#    no space, so it is not a whole-word boundary.)
$d.Content.Find.Execute("programas", $true, $false, $false, $false, $false,
                         $true, 1, $false, "programa", 2)
